$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1111.6666
$ws.Range("J43").Value = 1295
$ws.Range("L43").Value = 1295
$ws.Range("N43").Value = -1433

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6068.107
$ws.Range("J64").Value = 6380.8096
$ws.Range("L64").Value = 6380.8096
$ws.Range("N64").Value = -6876.8096

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6068.107
$ws.Range("J67").Value = 6380.8096
$ws.Range("L67").Value = 6380.8096
$ws.Range("N67").Value = -8096.8096

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5421.048
$ws.Range("I132").Value = 6534.893
$ws.Range("K132").Value = 19604.679
$ws.Range("M132").Value = -17074.679

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3204.5076
$ws.Range("I138").Value = 1732.7693
$ws.Range("J138").Value = 5254.4287
$ws.Range("K138").Value = 5198.3079
$ws.Range("L138").Value = 15763.2861
$ws.Range("M138").Value = -58.30789999999979
$ws.Range("N138").Value = -26043.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1677.5555
$ws.Range("I22").Value = 2066.5
$ws.Range("J22").Value = 899.6667
$ws.Range("K22").Value = 2066.5
$ws.Range("L22").Value = 899.6667
$ws.Range("M22").Value = -1767.5
$ws.Range("N22").Value = -1497.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9185838
$ws.Range("I32").Value = 4445335
$ws.Range("K32").Value = 4445335
$ws.Range("M32").Value = -4445048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2967.5454
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 99989
$ws.Range("J128").Value = 99989
$ws.Range("L128").Value = 99989
$ws.Range("N128").Value = -109949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3358.9832
$ws.Range("I132").Value = 2624.138
$ws.Range("K132").Value = 7872.414
$ws.Range("M132").Value = -5342.414

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2967.5454
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 29992
$ws.Range("J53").Value = 29992
$ws.Range("L53").Value = 29992
$ws.Range("N53").Value = -31140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11244865
$ws.Range("I134").Value = 2382155
$ws.Range("J134").Value = 55558412
$ws.Range("K134").Value = 7146465
$ws.Range("L134").Value = 166675236
$ws.Range("M134").Value = -7143930
$ws.Range("N134").Value = -166680306

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 858.9375
$ws.Range("I22").Value = 957.2857
$ws.Range("K22").Value = 957.2857
$ws.Range("M22").Value = -607.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7495.1333
$ws.Range("J31").Value = 11758.85
$ws.Range("L31").Value = 11758.85
$ws.Range("N31").Value = -12348.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7495.1333
$ws.Range("J34").Value = 11758.85
$ws.Range("L34").Value = 11758.85
$ws.Range("N34").Value = -12162.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 97643.664
$ws.Range("J75").Value = 113966.5
$ws.Range("L75").Value = 113966.5
$ws.Range("N75").Value = -115962.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 97643.664
$ws.Range("J78").Value = 113966.5
$ws.Range("L78").Value = 341899.5
$ws.Range("N78").Value = -351883.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 30468.334
$ws.Range("J124").Value = 29026.875
$ws.Range("L124").Value = 29026.875
$ws.Range("N124").Value = -33936.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1196.3549
$ws.Range("I132").Value = 1004.13043
$ws.Range("J132").Value = 1749
$ws.Range("K132").Value = 3012.39129
$ws.Range("L132").Value = 5247
$ws.Range("M132").Value = -482.39129
$ws.Range("N132").Value = -10307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2573.2222
$ws.Range("I134").Value = 2151.261
$ws.Range("K134").Value = 6453.782999999999
$ws.Range("M134").Value = -3918.782999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 953.1429000000001
$ws.Range("J107").Value = 954.86664
$ws.Range("L107").Value = 2864.59992
$ws.Range("N107").Value = -6704.59992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2001.7778
$ws.Range("J122").Value = 1883.8
$ws.Range("L122").Value = 16954.2
$ws.Range("N122").Value = -21854.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2455.818
$ws.Range("I132").Value = 2015.5714
$ws.Range("K132").Value = 18140.1426
$ws.Range("M132").Value = -15610.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1251.0333
$ws.Range("I102").Value = 1223.5927
$ws.Range("J102").Value = 1498
$ws.Range("K102").Value = 1223.5927
$ws.Range("L102").Value = 1498
$ws.Range("M102").Value = 398.4073000000001
$ws.Range("N102").Value = -4742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 29888
$ws.Range("J123").Value = 29888
$ws.Range("L123").Value = 29888
$ws.Range("N123").Value = -34788

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 152820
$ws.Range("J128").Value = 152820
$ws.Range("L128").Value = 152820
$ws.Range("N128").Value = -162780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2594.6
$ws.Range("I132").Value = 2225.389
$ws.Range("K132").Value = 6676.167
$ws.Range("M132").Value = -4146.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 745.2917
$ws.Range("I55").Value = 762.9286
$ws.Range("J55").Value = 720.6
$ws.Range("K55").Value = 762.9286
$ws.Range("L55").Value = 720.6
$ws.Range("M55").Value = -589.9286
$ws.Range("N55").Value = -1066.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3738.5
$ws.Range("I61").Value = 4926.5
$ws.Range("J61").Value = 1758.5
$ws.Range("K61").Value = 4926.5
$ws.Range("L61").Value = 1758.5
$ws.Range("M61").Value = -4724.5
$ws.Range("N61").Value = -2162.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3738.5
$ws.Range("I113").Value = 4926.5
$ws.Range("J113").Value = 1758.5
$ws.Range("K113").Value = 4926.5
$ws.Range("L113").Value = 1758.5
$ws.Range("M113").Value = -2756.5
$ws.Range("N113").Value = -6098.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2468531.8
$ws.Range("I122").Value = 3837728.2
$ws.Range("J122").Value = 3978
$ws.Range("K122").Value = 11513184.6
$ws.Range("L122").Value = 11934
$ws.Range("M122").Value = -11510734.6
$ws.Range("N122").Value = -16834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 89330
$ws.Range("J128").Value = 89330
$ws.Range("L128").Value = 89330
$ws.Range("N128").Value = -99290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3360.739
$ws.Range("I132").Value = 2099.8
$ws.Range("J132").Value = 4330.6924
$ws.Range("K132").Value = 6299.400000000001
$ws.Range("L132").Value = 12992.0772
$ws.Range("M132").Value = -3769.400000000001
$ws.Range("N132").Value = -18052.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 6262.5
$ws.Range("I40").Value = 5025
$ws.Range("K40").Value = 5025
$ws.Range("M40").Value = -4876

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1767.6923
$ws.Range("I122").Value = 1835.3182
$ws.Range("K122").Value = 5505.9546
$ws.Range("M122").Value = -3055.9546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 197597.25
$ws.Range("J128").Value = 197597.25
$ws.Range("L128").Value = 197597.25
$ws.Range("N128").Value = -207557.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 148997
$ws.Range("J130").Value = 148997
$ws.Range("L130").Value = 148997
$ws.Range("N130").Value = -159037

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3999.0334
$ws.Range("I132").Value = 3108.5908
$ws.Range("J132").Value = 6447.75
$ws.Range("K132").Value = 9325.7724
$ws.Range("L132").Value = 19343.25
$ws.Range("M132").Value = -6795.7724
$ws.Range("N132").Value = -24403.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 32992.473
$ws.Range("I136").Value = 2027.2963
$ws.Range("J136").Value = 125888
$ws.Range("K136").Value = 6081.8889
$ws.Range("L136").Value = 377664
$ws.Range("M136").Value = -3531.8889
$ws.Range("N136").Value = -382764
